$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style 9 source: Q4 (year header style)
# Style 11 source: A5 (bold sz9 text style, no borders, no numfmt)
# Style 10 source: A8 (plain sz9 text style, no borders, no numfmt)
# Style 17 (new font) built once from A8 + explicit theme color, reused afterward
# Style 18 (new xf) built once from A43 + right alignment, reused afterward

$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null
$ws.Range("R4").Value = 2021

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null
$ws.Range("R5").Value = 5.3

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R6").PasteSpecial(-4122) | Out-Null
$ws.Range("R6").Font.ThemeColor = 1
$ws.Range("R6").Value = 6.3

$ws.Range("R6").Copy() | Out-Null
$ws.Range("R7").PasteSpecial(-4122) | Out-Null
$ws.Range("R7").Value = 4.7

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R8").PasteSpecial(-4122) | Out-Null

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R9").PasteSpecial(-4122) | Out-Null
$ws.Range("R9").Value = 6.6

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R10").PasteSpecial(-4122) | Out-Null
$ws.Range("R10").Value = 7.5

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R11").PasteSpecial(-4122) | Out-Null
$ws.Range("R11").Value = 6.2

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R12").PasteSpecial(-4122) | Out-Null
$ws.Range("R12").Value = 11.8

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R13").PasteSpecial(-4122) | Out-Null
$ws.Range("R13").Value = 15.5

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R14").PasteSpecial(-4122) | Out-Null
$ws.Range("R14").Value = 9.6999999999999993

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R15").PasteSpecial(-4122) | Out-Null
$ws.Range("R15").Value = 6.3

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R16").PasteSpecial(-4122) | Out-Null
$ws.Range("R16").Value = 7.5

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").Value = 5.6

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R18").PasteSpecial(-4122) | Out-Null
$ws.Range("R18").Value = 6.3

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R19").PasteSpecial(-4122) | Out-Null
$ws.Range("R19").Value = 10.8

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R20").PasteSpecial(-4122) | Out-Null
$ws.Range("R20").Value = 4.3

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R21").PasteSpecial(-4122) | Out-Null
$ws.Range("R21").Value = 1.9

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R22").PasteSpecial(-4122) | Out-Null
$ws.Range("R22").Value = 3.1

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R23").PasteSpecial(-4122) | Out-Null
$ws.Range("R23").Value = 1.1000000000000001

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R24").PasteSpecial(-4122) | Out-Null
$ws.Range("R24").Value = 2.6

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R25").PasteSpecial(-4122) | Out-Null
$ws.Range("R25").Value = 3.8

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R26").PasteSpecial(-4122) | Out-Null
$ws.Range("R26").Value = 1.7

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R27").PasteSpecial(-4122) | Out-Null
$ws.Range("R27").Value = 5.3

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R28").PasteSpecial(-4122) | Out-Null
$ws.Range("R28").Value = 6.2

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R29").PasteSpecial(-4122) | Out-Null
$ws.Range("R29").Value = 4.8

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R30").PasteSpecial(-4122) | Out-Null
$ws.Range("R30").Value = 4.0999999999999996

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R31").PasteSpecial(-4122) | Out-Null
$ws.Range("R31").Value = 3.3

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R32").PasteSpecial(-4122) | Out-Null
$ws.Range("R32").Value = 4.9000000000000004

$ws.Range("A5").Copy() | Out-Null
$ws.Range("R33").PasteSpecial(-4122) | Out-Null
$ws.Range("R33").Value = 2.8

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R34").PasteSpecial(-4122) | Out-Null
$ws.Range("R34").Value = 3.4

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R35").PasteSpecial(-4122) | Out-Null
$ws.Range("R35").Value = 2.6

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R36").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R37").PasteSpecial(-4122) | Out-Null
$ws.Range("R37").Value = 15.7

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R38").PasteSpecial(-4122) | Out-Null
$ws.Range("R38").Value = 7.9

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R39").PasteSpecial(-4122) | Out-Null
$ws.Range("R39").Value = 4.5

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R40").PasteSpecial(-4122) | Out-Null
$ws.Range("R40").Value = 4.4000000000000004

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R41").PasteSpecial(-4122) | Out-Null
$ws.Range("R41").Value = 2.9

$ws.Range("A8").Copy() | Out-Null
$ws.Range("R42").PasteSpecial(-4122) | Out-Null
$ws.Range("R42").Value = 1.4

$ws.Range("A43").Copy() | Out-Null
$ws.Range("R43").PasteSpecial(-4122) | Out-Null
$ws.Range("R43").HorizontalAlignment = -4152
$ws.Range("R43").Value = "…"

$ws.Range("S1").Select() | Out-Null